$wb = $excel.ActiveWorkbook

# Update both the "展览" sheet and the "全部类型" sheet, which contain
# duplicated data rows. The "想去人数" (want-to-go count) column F changes:
#   Row 2: 129 -> 130
#   Row 3: 18  -> 19
#   Row 5: 22  -> 24
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 130
    $ws.Range("F3").Value = 19
    $ws.Range("F5").Value = 24
}
